$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "NBA, Wednesday 7th Feb 2024"

$ws.Range("A2").Value = "Toronto Raptors (17-33) vs Charlotte Hornets (10-39)"
$ws.Range("A3").Value = "Cleveland Cavaliers (32-16) vs Washington Wizards (9-40)"
$ws.Range("A4").Value = "Atlanta Hawks (22-28) vs Boston Celtics (38-12)"
$ws.Range("A5").Value = "San Antonio Spurs (10-40) vs Miami Heat (27-24)"
$ws.Range("A6").Value = "Golden State Warriors (22-25) vs Philadelphia 76ers (30-19)"
$ws.Range("A7").Value = "New Orleans Pelicans (29-21) vs Los Angeles Clippers (34-15)"
$ws.Range("A8").Value = "Detroit Pistons (6-43) vs Sacramento Kings (29-20)"

$ws.Range("B2").Value = "Toronto Raptors (43.48%)"
$ws.Range("B3").Value = "Cleveland Cavaliers (69.23%)"
$ws.Range("B4").Value = "Boston Celtics (88.46%)"
$ws.Range("B5").Value = "Miami Heat (53.85%)"
$ws.Range("B6").Value = "Philadelphia 76ers (68.00%)"
$ws.Range("B7").Value = "Los Angeles Clippers (82.61%)"
$ws.Range("B8").Value = "Sacramento Kings (63.64%)    "

$ws.Range("C2").Value = "Toronto Raptors (75.9%)"
$ws.Range("C3").Value = "Cleveland Cavaliers (87.8%)"
$ws.Range("C4").Value = "Boston Celtics (87.7%)"
$ws.Range("C5").Value = "Miami Heat (71.6%)"
$ws.Range("C6").Value = "Philadelphia 76ers (53.4%)"
$ws.Range("C7").Value = "Los Angeles Clippers (73.3%)"
$ws.Range("C8").Value = "Sacramento Kings (78.4%)    "

$ws.Range("C8").ClearFormats()

$ws.Range("C9").Select()

$wb.Save()
